$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format so numeric-looking strings (e.g. "1.000", "13.00")
# are preserved exactly as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '23.692.68'
$ws.Range("E2").Value = '  +1.98%  '

$ws.Range("D3").Value = '1.650.14'
$ws.Range("E3").Value = '  +2.83%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("D6").Value = '306.45'
$ws.Range("E6").Value = '  +0.64%  '

$ws.Range("D7").Value = '0.3782'
$ws.Range("E7").Value = '  +0.61%  '

$ws.Range("D8").Value = '52.99'
$ws.Range("E8").Value = '  +1.06%  '

$ws.Range("D9").Value = '0.3682'
$ws.Range("E9").Value = '  +1.61%  '

$ws.Range("D10").Value = '1.272'
$ws.Range("E10").Value = '  -0.18%  '

$ws.Range("D11").Value = '0.08169'
$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("D12").Value = '1.004'
$ws.Range("E12").Value = '  +0.34%  '

$ws.Range("D13").Value = '23.18'
$ws.Range("E13").Value = '  +1.13%  '

$ws.Range("D14").Value = '6.725'
$ws.Range("E14").Value = '  +1.83%  '

$ws.Range("D15").Value = '0.00001275'
$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("D16").Value = '7.407'
$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("D17").Value = '1.647.09'
$ws.Range("E17").Value = '  +2.63%  '

$ws.Range("D18").Value = '95.20'
$ws.Range("E18").Value = '  +1.29%  '

$ws.Range("D19").Value = '0.06937'
$ws.Range("E19").Value = '  +0.28%  '

$ws.Range("D20").Value = '18.41'
$ws.Range("E20").Value = '  +1.31%  '

$ws.Range("D21").Value = '6.592'
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").Value = '23.621.56'
$ws.Range("E23").Value = '  +1.71%  '

$ws.Range("D24").Value = '13.00'
$ws.Range("E24").Value = '  +0.57%  '

$ws.Range("D25").Value = '3.267'
$ws.Range("E25").Value = '  +6.45%  '

$ws.Range("D26").Value = '2.436'
$ws.Range("E26").Value = '  -0.53%  '

$ws.Range("D27").Value = '21.43'
$ws.Range("E27").Value = '  +1.26%  '

$ws.Range("D28").Value = '151.90'
$ws.Range("E28").Value = '  +1.04%  '

$ws.Range("D29").Value = '5.323'
$ws.Range("E29").Value = '  +0.76%  '

$ws.Range("D30").Value = '137.58'
$ws.Range("E30").Value = '  +1.76%  '

$ws.Range("B31").Value = 'WEMIXTOKEN'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").Value = '2.319'
$ws.Range("E31").Value = '  -3.02%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '7.053'
$ws.Range("E32").Value = '  +4.70%  '

$ws.Range("D33").Value = '1.831.18'
$ws.Range("E33").Value = '  +2.91%  '

$ws.Range("D34").Value = '11.04'
$ws.Range("E34").Value = '  +6.14%  '

$ws.Range("D35").Value = '0.9783'
$ws.Range("E35").Value = '  +1.62%  '

$ws.Range("D36").Value = '0.02899'
$ws.Range("E36").Value = '  +4.67%  '

$ws.Range("D37").Value = '6.388'
$ws.Range("E37").Value = '  +4.39%  '

$ws.Range("D38").Value = '0.2581'
$ws.Range("E38").Value = '  +2.35%  '

$ws.Range("D39").Value = '0.07350'
$ws.Range("E39").Value = '  -1.95%  '

$ws.Range("D40").Value = '0.08893'
$ws.Range("E40").Value = '  +0.99%  '

$ws.Range("D41").Value = '1.387'
$ws.Range("E41").Value = '  -1.49%  '

$ws.Range("D42").Value = '0.7213'
$ws.Range("E42").Value = '  +1.71%  '

$ws.Range("D43").Value = '12.79'
$ws.Range("E43").Value = '  +2.70%  '

$ws.Range("D44").Value = '16.55'
$ws.Range("E44").Value = '  +3.81%  '

$ws.Range("D45").Value = '0.6663'
$ws.Range("E45").Value = '  +1.90%  '

$ws.Range("D46").Value = '2.388'
$ws.Range("E46").Value = '  +2.39%  '

$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").Value = '4.028'
$ws.Range("E48").Value = '  +0.45%  '

$ws.Range("D49").Value = '0.08071'
$ws.Range("E49").Value = '  +1.59%  '

$ws.Range("D50").Value = '1.232'
$ws.Range("E50").Value = '  +2.15%  '

$ws.Range("E51").Value = '  -3.76%  '
